# CHG: updated risk plan with date of last change
#
# Add an "Aktualisierung:" (update) label with the date/sprint of the last
# change next to the sheet title on Tabelle1 (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Label cell
$ws.Range("B1").Value = "Aktualisierung:"

# Value cell - text that looks like a date stamp; give it a date-like
# number format (maps to built-in format id 14) even though the stored
# value stays a text string.
$ws.Range("C1").Value = "Sprint 11: 13.06.2013"
$ws.Range("C1").NumberFormat = "mm-dd-yy"

# Leave the freshly edited cell selected, as the author did.
$ws.Range("C1").Select()
